{"js": "// Replace the 25 three-digit-divided-by-one-digit division prompts in the\n// table with their new values (text before the trailing \"=\" changes; the\n// \"\u00f7\" and \"=\" characters stay put as part of the matched/replaced text).\nconst replacements = [\n  [\"140\u00f75=\", \"666\u00f76=\"],\n  [\"628\u00f74=\", \"543\u00f76=\"],\n  [\"987\u00f73=\", \"301\u00f78=\"],\n  [\"812\u00f79=\", \"549\u00f79=\"],\n  [\"478\u00f73=\", \"506\u00f77=\"],\n  [\"123\u00f78=\", \"953\u00f78=\"],\n  [\"264\u00f78=\", \"672\u00f77=\"],\n  [\"524\u00f76=\", \"962\u00f75=\"],\n  [\"218\u00f75=\", \"307\u00f75=\"],\n  [\"434\u00f79=\", \"858\u00f73=\"],\n  [\"955\u00f75=\", \"591\u00f78=\"],\n  [\"653\u00f79=\", \"921\u00f72=\"],\n  [\"918\u00f75=\", \"656\u00f77=\"],\n  [\"428\u00f78=\", \"711\u00f75=\"],\n  [\"290\u00f72=\", \"727\u00f77=\"],\n  [\"686\u00f77=\", \"595\u00f78=\"],\n  [\"904\u00f74=\", \"174\u00f79=\"],\n  [\"883\u00f72=\", \"879\u00f75=\"],\n  [\"509\u00f78=\", \"598\u00f76=\"],\n  [\"607\u00f74=\", \"645\u00f72=\"],\n  [\"890\u00f74=\", \"625\u00f74=\"],\n  [\"873\u00f72=\", \"273\u00f76=\"],\n  [\"840\u00f75=\", \"578\u00f78=\"],\n  [\"291\u00f78=\", \"311\u00f73=\"],\n  [\"565\u00f76=\", \"828\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 three-digit-divided-by-one-digit division prompts in the\n# table with their new values (text before the trailing \"=\" changes; the\n# \"\u00f7\" and \"=\" characters stay put as part of the matched/replaced text).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"140\u00f75=\"; New = \"666\u00f76=\" }\n    @{ Old = \"628\u00f74=\"; New = \"543\u00f76=\" }\n    @{ Old = \"987\u00f73=\"; New = \"301\u00f78=\" }\n    @{ Old = \"812\u00f79=\"; New = \"549\u00f79=\" }\n    @{ Old = \"478\u00f73=\"; New = \"506\u00f77=\" }\n    @{ Old = \"123\u00f78=\"; New = \"953\u00f78=\" }\n    @{ Old = \"264\u00f78=\"; New = \"672\u00f77=\" }\n    @{ Old = \"524\u00f76=\"; New = \"962\u00f75=\" }\n    @{ Old = \"218\u00f75=\"; New = \"307\u00f75=\" }\n    @{ Old = \"434\u00f79=\"; New = \"858\u00f73=\" }\n    @{ Old = \"955\u00f75=\"; New = \"591\u00f78=\" }\n    @{ Old = \"653\u00f79=\"; New = \"921\u00f72=\" }\n    @{ Old = \"918\u00f75=\"; New = \"656\u00f77=\" }\n    @{ Old = \"428\u00f78=\"; New = \"711\u00f75=\" }\n    @{ Old = \"290\u00f72=\"; New = \"727\u00f77=\" }\n    @{ Old = \"686\u00f77=\"; New = \"595\u00f78=\" }\n    @{ Old = \"904\u00f74=\"; New = \"174\u00f79=\" }\n    @{ Old = \"883\u00f72=\"; New = \"879\u00f75=\" }\n    @{ Old = \"509\u00f78=\"; New = \"598\u00f76=\" }\n    @{ Old = \"607\u00f74=\"; New = \"645\u00f72=\" }\n    @{ Old = \"890\u00f74=\"; New = \"625\u00f74=\" }\n    @{ Old = \"873\u00f72=\"; New = \"273\u00f76=\" }\n    @{ Old = \"840\u00f75=\"; New = \"578\u00f78=\" }\n    @{ Old = \"291\u00f78=\"; New = \"311\u00f73=\" }\n    @{ Old = \"565\u00f76=\"; New = \"828\u00f78=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $true, $pair.New, 2)\n}\n"}
